# Delete row 18 (the old "628994d9..." draft row) on the "Record Days" sheet,
# which shifts row 19 ("e432f0fe..." ready row) up to become the new row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record Days")

$ws.Rows.Item(18).Delete()
